$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Paragraph 39 ("SmartCard: Введение") -- restyle as Heading2, same text
# ---------------------------------------------------------------------------
$p39 = $d.Paragraphs(39)
$xml39 = '<w:p ' + $wNs + '>' +
  '<w:pPr>' +
    '<w:pStyle w:val="Heading2"/>' +
    '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
    '<w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/>' +
    '<w:jc w:val="center"/>' +
    '<w:textAlignment w:val="baseline"/>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="exo" w:hAnsi="exo" w:cs="open sans"/>' +
      '<w:caps/>' +
      '<w:color w:val="343434"/>' +
      '<w:spacing w:val="15"/>' +
      '<w:sz w:val="42"/>' +
      '<w:szCs w:val="42"/>' +
    '</w:rPr>' +
  '</w:pPr>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="exo" w:hAnsi="exo" w:cs="open sans"/>' +
      '<w:caps/>' +
      '<w:color w:val="343434"/>' +
      '<w:spacing w:val="15"/>' +
      '<w:sz w:val="42"/>' +
      '<w:szCs w:val="42"/>' +
    '</w:rPr>' +
    '<w:t xml:space="preserve">SmartCard: Введение</w:t>' +
  '</w:r>' +
  '</w:p>'
$p39.Range.InsertXML($xml39)

# ---------------------------------------------------------------------------
# 2) Paragraph 40 ("Мы подготовили ...") -- restyle as Heading2, new text
# ---------------------------------------------------------------------------
$p40 = $d.Paragraphs(40)
$xml40 = '<w:p ' + $wNs + '>' +
  '<w:pPr>' +
    '<w:pStyle w:val="Heading2"/>' +
    '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
    '<w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/>' +
    '<w:jc w:val="center"/>' +
    '<w:textAlignment w:val="baseline"/>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="exo" w:hAnsi="exo"/>' +
      '<w:caps/>' +
      '<w:color w:val="343434"/>' +
      '<w:spacing w:val="15"/>' +
      '<w:sz w:val="42"/>' +
      '<w:szCs w:val="42"/>' +
    '</w:rPr>' +
  '</w:pPr>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="exo" w:hAnsi="exo"/>' +
      '<w:caps/>' +
      '<w:color w:val="343434"/>' +
      '<w:spacing w:val="15"/>' +
      '<w:sz w:val="42"/>' +
      '<w:szCs w:val="42"/>' +
    '</w:rPr>' +
    '<w:t xml:space="preserve">SMARTCASH FOR BUSINESS</w:t>' +
  '</w:r>' +
  '</w:p>'
$p40.Range.InsertXML($xml40)

# ---------------------------------------------------------------------------
# 3) Insert 9 brand-new paragraphs right after paragraph 40, before the
#    pre-existing trailing empty paragraph (which must stay untouched).
#    The LAST paragraph supplied here is itself empty, which keeps the
#    runtime from merging it into -- and clobbering -- the old trailing
#    paragraph that follows.
# ---------------------------------------------------------------------------
$p40b = $d.Paragraphs(40)
$insertPoint = $d.Range($p40b.Range.End, $p40b.Range.End)

$newBlock =
  # empty centered Heading2 paragraph
  ('<w:p ' + $wNs + '>' +
    '<w:pPr>' +
      '<w:pStyle w:val="Heading2"/>' +
      '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
      '<w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/>' +
      '<w:jc w:val="center"/>' +
      '<w:textAlignment w:val="baseline"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="exo" w:hAnsi="exo" w:cs="open sans"/>' +
        '<w:caps/>' +
        '<w:color w:val="343434"/>' +
        '<w:spacing w:val="15"/>' +
        '<w:sz w:val="42"/>' +
        '<w:szCs w:val="42"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
  '</w:p>') +
  # "SMARTCARD TUTORIAL VIDEOS" -- Heading1, centered, sz 51
  ('<w:p ' + $wNs + '>' +
    '<w:pPr>' +
      '<w:pStyle w:val="Heading1"/>' +
      '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
      '<w:spacing w:before="0" w:beforeAutospacing="0" w:after="210" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/>' +
      '<w:jc w:val="center"/>' +
      '<w:textAlignment w:val="baseline"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="exo" w:hAnsi="exo"/>' +
        '<w:caps/>' +
        '<w:color w:val="343434"/>' +
        '<w:spacing w:val="15"/>' +
        '<w:sz w:val="51"/>' +
        '<w:szCs w:val="51"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="exo" w:hAnsi="exo"/>' +
        '<w:caps/>' +
        '<w:color w:val="343434"/>' +
        '<w:spacing w:val="15"/>' +
        '<w:sz w:val="51"/>' +
        '<w:szCs w:val="51"/>' +
      '</w:rPr>' +
      '<w:t xml:space="preserve">SMARTCARD TUTORIAL VIDEOS</w:t>' +
    '</w:r>' +
  '</w:p>') +
  # "DOWNLOAD SMARTPAY APP" -- Heading1, not centered, sz 51
  ('<w:p ' + $wNs + '>' +
    '<w:pPr>' +
      '<w:pStyle w:val="Heading1"/>' +
      '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
      '<w:spacing w:before="0" w:beforeAutospacing="0" w:after="210" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/>' +
      '<w:textAlignment w:val="baseline"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="exo" w:hAnsi="exo"/>' +
        '<w:caps/>' +
        '<w:color w:val="343434"/>' +
        '<w:spacing w:val="15"/>' +
        '<w:sz w:val="51"/>' +
        '<w:szCs w:val="51"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="exo" w:hAnsi="exo"/>' +
        '<w:caps/>' +
        '<w:color w:val="343434"/>' +
        '<w:spacing w:val="15"/>' +
        '<w:sz w:val="51"/>' +
        '<w:szCs w:val="51"/>' +
      '</w:rPr>' +
      '<w:t xml:space="preserve">DOWNLOAD SMARTPAY APP</w:t>' +
    '</w:r>' +
  '</w:p>') +
  # "Accept SmartCash ..." -- NormalWeb
  ('<w:p ' + $wNs + '>' +
    '<w:pPr>' +
      '<w:pStyle w:val="NormalWeb"/>' +
      '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
      '<w:spacing w:before="204" w:beforeAutospacing="0" w:after="204" w:afterAutospacing="0"/>' +
      '<w:textAlignment w:val="baseline"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="open sans" w:hAnsi="open sans" w:cs="open sans"/>' +
        '<w:color w:val="252525"/>' +
        '<w:sz w:val="21"/>' +
        '<w:szCs w:val="21"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="open sans" w:hAnsi="open sans" w:cs="open sans"/>' +
        '<w:color w:val="252525"/>' +
        '<w:sz w:val="21"/>' +
        '<w:szCs w:val="21"/>' +
      '</w:rPr>' +
      '<w:t xml:space="preserve">Accept SmartCash as a payment option in your business with zero fees using a simple SmartPay app. The SmartPay app is available for use anywhere in the world, all it requires is access to the internet.</w:t>' +
    '</w:r>' +
  '</w:p>') +
  # "SMARTCARD SHOP" -- Heading2, not centered, sz 42
  ('<w:p ' + $wNs + '>' +
    '<w:pPr>' +
      '<w:pStyle w:val="Heading2"/>' +
      '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
      '<w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/>' +
      '<w:textAlignment w:val="baseline"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="exo" w:hAnsi="exo"/>' +
        '<w:caps/>' +
        '<w:color w:val="343434"/>' +
        '<w:spacing w:val="15"/>' +
        '<w:sz w:val="42"/>' +
        '<w:szCs w:val="42"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="exo" w:hAnsi="exo"/>' +
        '<w:caps/>' +
        '<w:color w:val="343434"/>' +
        '<w:spacing w:val="15"/>' +
        '<w:sz w:val="42"/>' +
        '<w:szCs w:val="42"/>' +
      '</w:rPr>' +
      '<w:t xml:space="preserve">SMARTCARD SHOP</w:t>' +
    '</w:r>' +
  '</w:p>') +
  # "An online store ..." -- NormalWeb
  ('<w:p ' + $wNs + '>' +
    '<w:pPr>' +
      '<w:pStyle w:val="NormalWeb"/>' +
      '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
      '<w:spacing w:before="204" w:beforeAutospacing="0" w:after="204" w:afterAutospacing="0"/>' +
      '<w:textAlignment w:val="baseline"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="open sans" w:hAnsi="open sans" w:cs="open sans"/>' +
        '<w:color w:val="252525"/>' +
        '<w:sz w:val="21"/>' +
        '<w:szCs w:val="21"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="open sans" w:hAnsi="open sans" w:cs="open sans"/>' +
        '<w:color w:val="252525"/>' +
        '<w:sz w:val="21"/>' +
        '<w:szCs w:val="21"/>' +
      '</w:rPr>' +
      '<w:t xml:space="preserve">An online store where you can purchase physical SmartCards.</w:t>' +
    '</w:r>' +
  '</w:p>') +
  # "READY TO GET STARTED?" -- Heading2, not centered, sz 42
  ('<w:p ' + $wNs + '>' +
    '<w:pPr>' +
      '<w:pStyle w:val="Heading2"/>' +
      '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
      '<w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/>' +
      '<w:textAlignment w:val="baseline"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="exo" w:hAnsi="exo"/>' +
        '<w:caps/>' +
        '<w:color w:val="343434"/>' +
        '<w:spacing w:val="15"/>' +
        '<w:sz w:val="42"/>' +
        '<w:szCs w:val="42"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="exo" w:hAnsi="exo"/>' +
        '<w:caps/>' +
        '<w:color w:val="343434"/>' +
        '<w:spacing w:val="15"/>' +
        '<w:sz w:val="42"/>' +
        '<w:szCs w:val="42"/>' +
      '</w:rPr>' +
      '<w:t xml:space="preserve">READY TO GET STARTED?</w:t>' +
    '</w:r>' +
  '</w:p>') +
  # empty centered Heading2 paragraph
  ('<w:p ' + $wNs + '>' +
    '<w:pPr>' +
      '<w:pStyle w:val="Heading2"/>' +
      '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
      '<w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/>' +
      '<w:jc w:val="center"/>' +
      '<w:textAlignment w:val="baseline"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="exo" w:hAnsi="exo" w:cs="open sans"/>' +
        '<w:caps/>' +
        '<w:color w:val="343434"/>' +
        '<w:spacing w:val="15"/>' +
        '<w:sz w:val="42"/>' +
        '<w:szCs w:val="42"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
  '</w:p>') +
  # final empty "inherit" paragraph (stays empty so it does not merge into
  # -- and overwrite -- the pre-existing trailing paragraph that follows it)
  ('<w:p ' + $wNs + '>' +
    '<w:pPr>' +
      '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
      '<w:spacing w:line="396" w:lineRule="atLeast"/>' +
      '<w:textAlignment w:val="baseline"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="open sans"/>' +
        '<w:color w:val="252525"/>' +
        '<w:sz w:val="21"/>' +
        '<w:szCs w:val="21"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
  '</w:p>')

$insertPoint.InsertXML($newBlock)
